# Replace the M2Doc field ( { m:x } ) which is currently encoded as a real
# Word field (fldChar begin / instrText runs / fldChar end) with plain
# literal-text runs spelling out the same token: "{", "m", ":x", "}".
#
# This mirrors the TokenIteratorFieldRewriterSplit change: the template no
# longer relies on a Word field, the braces/tag are stored as plain w:t
# runs instead.

$d = $word.ActiveDocument

# Find the paragraph that hosts the field (there is exactly one field in
# this template, but walk the collection instead of hardcoding an index so
# the script keeps working if surrounding paragraphs shift).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
    }
}

$paraRange = $targetPara.Range
# Exclude the trailing paragraph mark so we only rewrite the paragraph's
# content (the field + its fldChar begin/end + instrText runs) and keep
# the paragraph itself intact.
$contentRange = $d.Range($paraRange.Start, $paraRange.End - 1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    '<w:r><w:t>:x</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
    '</w:p>'

$contentRange.InsertXML($newParaXml)
